# Generate Report for Archive
#
# Updates the "Ready for handoff" status to "In Translation" across the
# Overview / zh-cn / de-de sheets, and shrinks the corresponding status
# column widths to match the regenerated report layout.

$wb  = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update status text (all cells sharing the "Ready for handoff" string) ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# --- Narrow the Status / locale columns ---
# Target stored column width is ~13.41 characters; the host subtracts a
# fixed 5/6-character padding when persisting ColumnWidth, so we offset
# the value we assign accordingly.
$newColumnWidth = 13.4101845877511 - (5/6)

$wsOverview.Range("E1").ColumnWidth = $newColumnWidth
$wsOverview.Range("F1").ColumnWidth = $newColumnWidth
$wsZhCn.Range("C1").ColumnWidth     = $newColumnWidth
$wsDeDe.Range("C1").ColumnWidth     = $newColumnWidth
